$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Add new row of data (row 18) for "Remove Duplicates from Sorted List"
$ws1.Range("B18").Value = "Remove Duplicates from Sorted List"
$ws1.Range("C18").Value = 1
$ws1.Range("D18").Value = 53
$ws1.Range("E18").Value = 0.83
$ws1.Range("F18").Value = 16.39
$ws1.Range("G18").Value = 0.61
$ws1.Range("H18").Value = "https://leetcode.com/problems/remove-duplicates-from-sorted-list/submissions/"

# Update the selection on Sheet1 to H21
$ws1.Range("H21").Select()
